# Revert "Powerpoint writer: consolidate text run nodes."
#
# Splits previously-merged text runs back into separate <a:r> nodes so
# that a trailing space is its own run rather than being glued onto the
# preceding word. Re-assigning Characters(start, length).Text with the
# same text forces the writer to emit a dedicated run for that span.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: "A " + "slide"  ->  "A" + " " + "slide" ---
$title = $s.Shapes.Item(1)
$trTitle = $title.TextFrame.TextRange
$trTitle.Characters(1, 1).Text = "A"
$trTitle.Characters(2, 1).Text = " "

# --- Caption textbox: "Followed " + "by " + "a " + "picture" ---
#     ->  "Followed" + " " + "by" + " " + "a" + " " + "picture"
$caption = $s.Shapes.Item(4)
$trCaption = $caption.TextFrame.TextRange
$trCaption.Characters(1, 8).Text = "Followed"
$trCaption.Characters(9, 1).Text = " "
$trCaption.Characters(10, 2).Text = "by"
$trCaption.Characters(12, 1).Text = " "
$trCaption.Characters(13, 1).Text = "a"
$trCaption.Characters(14, 1).Text = " "
$trCaption.Characters(15, 7).Text = "picture"
